$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 131139047
$ws.Range("B13").Value = 57881
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 100049
$ws.Range("F13").Value = "Spillkråka"
$ws.Range("G13").Value = "Dryocopus martius"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("I13").Value = "'1"
$ws.Range("M13").Value = "födosökande"
$ws.Range("P13").Value = "Korshamnsgrundet, Korshamnsgrundet, Sm"
$ws.Range("Q13").Value = 592528
$ws.Range("R13").Value = 6320591
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = "Kalmar"
$ws.Range("U13").Value = "Mönsterås"
$ws.Range("V13").Value = "Småland"
$ws.Range("W13").Value = "Mönsterås"
$ws.Range("Y13").Value = "'2026-02-13"
$ws.Range("AA13").Value = "'2026-02-13"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AW13").Value = "Lars Engström"
$ws.Range("AX13").Value = "Lars Engström"

# Row 14
$ws.Range("A14").Value = 131146326
$ws.Range("B14").Value = 57830
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 100067
$ws.Range("F14").Value = "Havsörn"
$ws.Range("G14").Value = "Haliaeetus albicilla"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("I14").Value = "'2"
$ws.Range("K14").Value = "2K+"
$ws.Range("M14").Value = "lockläte, övriga läten"
$ws.Range("P14").Value = "Korshamn, Björnö, Mönsterås, Sm"
$ws.Range("Q14").Value = 592599
$ws.Range("R14").Value = 6320593
$ws.Range("S14").Value = 25
$ws.Range("T14").Value = "Kalmar"
$ws.Range("U14").Value = "Mönsterås"
$ws.Range("V14").Value = "Småland"
$ws.Range("W14").Value = "Mönsterås"
$ws.Range("Y14").Value = "'2026-02-13"
$ws.Range("Z14").Value = "13:00"
$ws.Range("AA14").Value = "'2026-02-13"
$ws.Range("AB14").Value = "15:30"
$ws.Range("AC14").Value = "Kom över oss och flög ut mot Gryssholm"
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = "Jan Brenander"
$ws.Range("AX14").Value = "Jan Brenander"

# Row 15
$ws.Range("A15").Value = 131146287
$ws.Range("B15").Value = 58043
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 103021
$ws.Range("F15").Value = "Talltita"
$ws.Range("G15").Value = "Poecile montanus"
$ws.Range("H15").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I15").Value = "'1"
$ws.Range("M15").Value = "lockläte, övriga läten"
$ws.Range("P15").Value = "Korshamn, Björnö, Mönsterås, Sm"
$ws.Range("Q15").Value = 592599
$ws.Range("R15").Value = 6320593
$ws.Range("S15").Value = 25
$ws.Range("T15").Value = "Kalmar"
$ws.Range("U15").Value = "Mönsterås"
$ws.Range("V15").Value = "Småland"
$ws.Range("W15").Value = "Mönsterås"
$ws.Range("Y15").Value = "'2026-02-13"
$ws.Range("Z15").Value = "13:00"
$ws.Range("AA15").Value = "'2026-02-13"
$ws.Range("AB15").Value = "15:30"
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AW15").Value = "Jan Brenander"
$ws.Range("AX15").Value = "Jan Brenander"

# Row 16
$ws.Range("A16").Value = 131146362
$ws.Range("B16").Value = 91829
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5442
$ws.Range("F16").Value = "Tallticka"
$ws.Range("G16").Value = "Porodaedalea pini"
$ws.Range("H16").Value = "(Brot.) Murrill"
$ws.Range("I16").Value = "'1"
$ws.Range("P16").Value = "Korshamn, Björnö, Mönsterås, Sm"
$ws.Range("Q16").Value = 592629
$ws.Range("R16").Value = 6320625
$ws.Range("S16").Value = 25
$ws.Range("T16").Value = "Kalmar"
$ws.Range("U16").Value = "Mönsterås"
$ws.Range("V16").Value = "Småland"
$ws.Range("W16").Value = "Mönsterås"
$ws.Range("Y16").Value = "'2026-02-13"
$ws.Range("AA16").Value = "'2026-02-13"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = "Jan Brenander"
$ws.Range("AX16").Value = "Jan Brenander"

